$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("G3").Value = 80
$ws.Range("F4").Value = 317
$ws.Range("F5").Value = 215
$ws.Range("F6").Value = 520
$ws.Range("F7").Value = 1294
$ws.Range("F8").Value = 609
$ws.Range("F9").Value = 328
$ws.Range("F11").Value = 146
$ws.Range("F12").Value = 393
$ws.Range("F13").Value = 5981
$ws.Range("F15").Value = 19
$ws.Range("F16").Value = 1856
$ws.Range("F17").Value = 4445
$ws.Range("F18").Value = 452
$ws.Range("F21").Value = 5183
$ws.Range("F22").Value = 6722
$ws.Range("F24").Value = 1071
$ws.Range("F25").Value = 728
$ws.Range("F26").Value = 3900
$ws.Range("F27").Value = 525
$ws.Range("F31").Value = 1024
$ws.Range("F32").Value = 1460
$ws.Range("F33").Value = 524
$ws.Range("F34").Value = 629
$ws.Range("F35").Value = 1648
$ws.Range("F37").Value = 1814
$ws.Range("F39").Value = 1194
$ws.Range("F43").Value = 3565
$ws.Range("F45").Value = 325
$ws.Range("F48").Value = 70
$ws.Range("F49").Value = 3923

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1240
$ws.Range("F9").Value = 19
$ws.Range("F27").Value = 82

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4199

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4199
$ws.Range("F3").Value = 2715
$ws.Range("G4").Value = 80
$ws.Range("F5").Value = 317
$ws.Range("F6").Value = 1240
$ws.Range("F8").Value = 215
$ws.Range("F9").Value = 520
$ws.Range("F11").Value = 1294
$ws.Range("F12").Value = 19
$ws.Range("F13").Value = 609
$ws.Range("F14").Value = 328
$ws.Range("F15").Value = 146
$ws.Range("F16").Value = 393
$ws.Range("F18").Value = 1856
$ws.Range("F19").Value = 4445
$ws.Range("F20").Value = 5183
$ws.Range("F21").Value = 5183
$ws.Range("F23").Value = 1071
$ws.Range("F24").Value = 728
$ws.Range("F25").Value = 3900
$ws.Range("F26").Value = 525
$ws.Range("F29").Value = 1024
$ws.Range("F30").Value = 1460
$ws.Range("F31").Value = 524
$ws.Range("F32").Value = 629
$ws.Range("F33").Value = 1648
$ws.Range("F35").Value = 1814
$ws.Range("F42").Value = 82
$ws.Range("F43").Value = 3565
$ws.Range("F48").Value = 70
$ws.Range("F50").Value = 3923
